$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the last record (row 8) into a new row 9, then fix up the
# phone number to the cleaned-up value for that contact.
$ws.Range("A8:K8").Copy($ws.Range("A9:K9"))
$ws.Range("I9").Value = "(73) 9 8824-8659"

# Widen the new "Vendedor" helper column (K) to fit its content.
$ws.Columns.Item(11).ColumnWidth = 30.83

# Highlight the header row so it stands out.
$ws.Range("A1:K1").Interior.Color = 65535

# Rescale the view and leave the selection parked a few rows below the data.
$excel.ActiveWindow.Zoom = 85
$ws.Range("B13").Select()
